# Fix circular-reference SUM ranges on the LOE sheet (Management rollup rows
# and the TOTAL HOURS / TOTAL COST row) so they sum the actual task rows
# (E3:E29 / E3:E30 / E3:E31) instead of ranges that looped back onto the
# formula cells themselves.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LOE")

$ws.Range("D30").Formula = '=ROUND(SUM(E3:E29)*$C$30,0)'
$ws.Range("E30").Formula = '=ROUND(SUM(E3:E29)*0.25,0)'

$ws.Range("D31").Formula = '=ROUND(SUM(E3:E30)*0.20,0)'
$ws.Range("E31").Formula = '=ROUND(SUM(E3:E30)*0.20,0)'

$ws.Range("E32").Formula = '=SUM(E3:E31)'
$ws.Range("G32").Formula = '=TEXT(SUM(G3:G31),"$#,##0")'
